$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: header labels ----
$ws.Range("A1").Value = "datnum"
$ws.Range("B1").Value = "datname"
$ws.Range("C1").Value = "time"
$ws.Range("D1").Value = "picklepath"
$ws.Range("E1").Value = "x_label"
$ws.Range("F1").Value = "y_label"
$ws.Range("G1").Value = "dim"
$ws.Range("H1").Value = "time_elapsed"

# ---- Row 2 ----
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "base"
$ws.Range("C2").Value = 1577779312.384152
$ws.Range("D2").Value = "pathtopickle"
$ws.Range("E2").Value = "xlabel"
$ws.Range("F2").Value = "ylabel"

# ---- Row 3 ----
$ws.Range("A3").Value = 2700
$ws.Range("B3").Value = "base"
$ws.Range("E3").Value = "FD_SDP/1000mV"
$ws.Range("F3").Value = "Repeats (mV)"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 13.421

# ---- Apply the same bold/bordered header style (style index 1, as
#      already used by B1:D1 in the source file) to every new header
#      cell in row 1 (A1, E1:H1), and to the two index columns (A & B)
#      for the data rows -- matching the two-level MultiIndex
#      (datnum, datname) that now labels every row. ----
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("E1:H1").PasteSpecial(-4122)
$ws.Range("A2:B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
